$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "freesearch"
$ws.Range("D42").Value = "[Microsoft Spy++]Spy++ 용도 및 사용방법"
$ws.Range("D51").Value = "[VSCODE] PDF 파일 보기, vscode-pdf"
$ws.Range("E51").Value = "https://bskyvision.com/entry/VSCODE-PDF-%ED%8C%8C%EC%9D%BC-%EB%B3%B4%EA%B8%B0-vscode-pdf"
